$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping-clase-de-propietario")

$ws.Range("A1").Value = "Una persona"
$ws.Range("A2").Value = "Un organismo público"
$ws.Range("A3").Value = "La comunidad"
$ws.Range("A4").Value = "Una sociedad"
